$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New detection-result columns (F/G) for the last three data rows ---
$ws.Range("F14").Value = 0.751
$ws.Range("G14").Value = 0.66

$ws.Range("F15").Value = 0.738
$ws.Range("G15").Value = 0.608

$ws.Range("F16").Value = 0.753
$ws.Range("G16").Value = 0.625

# --- Summary row: average of the last three rows per column ---
$ws.Range("B18").Formula = "=SUM(B14:B16)/3"
$ws.Range("C18:G18").Formula = "=SUM(C14:C16)/3"

# --- Move/resize the chart to its new anchor position ---
$co = $ws.ChartObjects().Item(1)
$co.Left = 503.28574803149604
$co.Top = 72.03566929133858
$co.Width = 385.0625
$co.Height = 216.10708661417323

# --- Update the active selection ---
$ws.Range("P17").Select()
